# Update sig_dates holidays: mark the good Friday (row 7) and
# Washington's birthday (row 13) rows as market-closed days (vol = -1),
# matching the sheet's "Input a -1 into vol to signify market closed day" convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = -1
$ws.Range("B13").Value = -1

# Leave the selection where the author last left it while working on BT.
$ws.Range("F11").Select()
